$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 270, shifting existing rows 270:295 down to 271:296
$ws.Range("A270").EntireRow.Insert()

# Populate the newly inserted row 270 with the new data record
$ws.Range("A270").Value = 10
$ws.Range("B270").Value = "Vega Modelo de Temuco"
$ws.Range("C270").Value = "La Araucanía"
$ws.Range("D270").Value = 44578
$ws.Range("E270").Value = 9
$ws.Range("F270").Value = 100112040
$ws.Range("G270").Value = "Cilantro"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 40
$ws.Range("K270").Value = 8000
$ws.Range("L270").Value = 8000
$ws.Range("M270").Value = 8000
$ws.Range("N270").Value = "$/docena de atados (2 kilos)"
$ws.Range("O270").Value = "Provincia de Cautín"
$ws.Range("P270").Value = 4000
$ws.Range("Q270").Value = 2
$ws.Range("R270").Value = "Hortaliza"
